# Experiment order generation script
# Regenerates the randomized per-task stim orders for participant_33 and
# renames each task-order tab to match the freshly generated run id.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (was GNG_TO-16512555771186125): now holds the vSAT task order.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vSAT_TO-16515890102094724"
$ws1.Range("B2").Value = "vSAT_stims-16515890101782236.csv"
$ws1.Range("B3").Value = "SAT_stims-16515890101470032.csv"
$ws1.Range("B4").Value = "SAT_stims-16515890101625996.csv"
$ws1.Range("B5").Value = "vSAT_stims-16515890101938486.csv"

# ---------------------------------------------------------------------
# Sheet 2 (was NB_TO-1651255579562095): now holds the GNG task order.
# Shrinks from 8 data rows to 4 data rows.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "GNG_TO-16515890102407246"
$ws2.Range("B2").Value = "go_stims-16515890102094724.csv"
$ws2.Range("B3").Value = "GNG_stims-16515890102251015.csv"
$ws2.Range("B4").Value = "go_stims-16515890102251015.csv"
$ws2.Range("B5").Value = "GNG_stims-16515890102407246.csv"
$ws2.Range("A6:B10").Clear()

# ---------------------------------------------------------------------
# Sheet 3 (was RS_TO-1651255579562095): resting-state order, the two
# conditions swap places.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16515890102407246"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# ---------------------------------------------------------------------
# Sheet 4 (was TOL_TO-16512555796083114): still the N-back style task
# order (MM/ZM stims), row count unchanged, just fresh filenames.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16515890102876"
$ws4.Range("B2").Value = "MM_stims-1651589010256349.csv"
$ws4.Range("B3").Value = "ZM_stims-16515890102407246.csv"
$ws4.Range("B4").Value = "MM_stims-16515890102719734.csv"
$ws4.Range("B5").Value = "ZM_stims-1651589010256349.csv"
$ws4.Range("B6").Value = "MM_stims-16515890102876.csv"
$ws4.Range("B7").Value = "ZM_stims-16515890102719734.csv"

# ---------------------------------------------------------------------
# Sheet 5 (was vSAT_TO-16512555796801126): now holds the NB task order.
# Grows from 4 data rows to 9 data rows.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "NB_TO-16515890116039078"
$ws5.Range("B2").Value = "TB-16515890114923708.csv"
$ws5.Range("B3").Value = "OB-16515890107305393.csv"
$ws5.Range("B4").Value = "ZB-match_9-165158901029055.csv"
$ws5.Range("B5").Value = "OB-16515890109982982.csv"

$ws5.Range("A2").Copy($ws5.Range("A6"))
$ws5.Range("A2").Copy($ws5.Range("A7"))
$ws5.Range("A2").Copy($ws5.Range("A8"))
$ws5.Range("A2").Copy($ws5.Range("A9"))
$ws5.Range("A2").Copy($ws5.Range("A10"))

$ws5.Range("A6").Value = 4
$ws5.Range("B6").Value = "TB-1651589011588244.csv"
$ws5.Range("A7").Value = 5
$ws5.Range("B7").Value = "ZB-match_7-1651589010472.csv"
$ws5.Range("A8").Value = 6
$ws5.Range("B8").Value = "TB-16515890115236182.csv"
$ws5.Range("A9").Value = 7
$ws5.Range("B9").Value = "ZB-match_4-16515890104407508.csv"
$ws5.Range("A10").Value = 8
$ws5.Range("B10").Value = "OB-16515890111723404.csv"
